$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete rows 5 and 6 (companies removed from dataset) ---
$ws.Range("A5:A6").EntireRow.Delete()

# --- Row 2 updates ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value = 0.22
$ws.Range("G2").Value = -0.1790464240903388
$ws.Range("H2").Value = -0.1838143036386449
$ws.Range("I2").Value = -0.2104767879548306
$ws.Range("J2").Value = -0.2104767879548306
$ws.Range("K2").Value = -5.649999999999999
$ws.Range("L2").Value = -0.1772271016311167
$ws.Range("U2").Value = 2.222
$ws.Range("V2").Value = 0.08693270735524257
$ws.Range("W2").Value = -46.26682724474775
$ws.Range("X2").Value = 0.07764728542239972
$ws.Range("Y2").Value = -46.34447453017015
$ws.Range("Z2").Value = 7.877440079070917
$ws.Range("AA2").Value = -3.090391251260817
$ws.Range("AB2").Value = 0.07209091411472041
$ws.Range("AC2").Value = -3.162482165375538
$ws.Range("AD2").Value = 2.94
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2.94
$ws.Range("AG2").Value = 0.718
$ws.Range("AH2").Value = 0.1031578947368421
$ws.Range("AI2").Value = 1.563829787234043
$ws.Range("AJ2").Value = 0.02732323616713601
$ws.Range("AK2").Value = -2.099415204678359
$ws.Range("AL2").Value = 0.205
$ws.Range("AM2").Value = 0.203
$ws.Range("AN2").Value = -0.5526315789473684
$ws.Range("AO2").Value = -32.73170731707317
$ws.Range("AP2").Value = -0.1349624060150376
$ws.Range("AQ2").Value = -33.05418719211822

# --- Row 3 updates ---
$ws.Range("B3").Value = "Y Ventures Group Ltd. (Catalist:1F1)"
$ws.Range("G3").Value = -0.06557823129251701
$ws.Range("H3").Value = -0.0707482993197279
$ws.Range("I3").Value = -0.0554421768707483
$ws.Range("J3").Value = -0.0554421768707483
$ws.Range("K3").Value = -0.76
$ws.Range("L3").Value = -0.02585034013605442
$ws.Range("U3").Value = 2.16
$ws.Range("V3").Value = 0.1375796178343949
$ws.Range("W3").Value = -0.2695035460992908
$ws.Range("X3").Value = 0.07528958334349715
$ws.Range("Y3").Value = -0.3447931294427879
$ws.Range("Z3").Value = 9.333333333333334
$ws.Range("AA3").Value = -0.5174603174603174
$ws.Range("AB3").Value = 0.07224225050105636
$ws.Range("AC3").Value = -0.5897025679613738
$ws.Range("AD3").Value = 1.25
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1.25
$ws.Range("AG3").Value = -0.9100000000000001
$ws.Range("AH3").Value = 0.07374631268436578
$ws.Range("AI3").Value = 0.2626050420168067
$ws.Range("AJ3").Value = -0.06152805949966195
$ws.Range("AK3").Value = -0.3500000000000001
$ws.Range("AL3").Value = 0.049
$ws.Range("AM3").Value = 0.047
$ws.Range("AN3").Value = -0.8116883116883117
$ws.Range("AO3").Value = -33.26530612244898
$ws.Range("AP3").Value = 0.5909090909090909
$ws.Range("AQ3").Value = -34.68085106382978

# --- Row 4 updates ---
$ws.Range("B4").Value = "LifeBrandz Ltd. (Catalist:1D3)"
$ws.Range("D4").Value = 0.22
$ws.Range("G4").Value = -1.524193548387097
$ws.Range("H4").Value = -1.524193548387097
$ws.Range("I4").Value = -2.048387096774194
$ws.Range("J4").Value = -2.048387096774194
$ws.Range("K4").Value = -4.89
$ws.Range("L4").Value = -1.971774193548387
$ws.Range("U4").Value = 0.062
$ws.Range("V4").Value = 0.006288032454361055
$ws.Range("W4").Value = -92.26415094339622
$ws.Range("X4").Value = 0.08000498750130229
$ws.Range("Y4").Value = -92.34415593089751
$ws.Range("Z4").Value = 2.764771460423634
$ws.Range("AA4").Value = -5.663322185061316
$ws.Range("AB4").Value = 0.07193957772838447
$ws.Range("AC4").Value = -5.735261762789701
$ws.Range("AD4").Value = 1.69
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 1.69
$ws.Range("AG4").Value = 1.628
$ws.Range("AH4").Value = 0.1463203463203463
$ws.Range("AI4").Value = -0.5868055555555555
$ws.Range("AJ4").Value = 0.1417130919220056
$ws.Range("AK4").Value = -0.5533650577838205
$ws.Range("AL4").Value = 0.156
$ws.Range("AM4").Value = 0.156
$ws.Range("AN4").Value = -0.4470899470899471
$ws.Range("AO4").Value = -32.56410256410256
$ws.Range("AP4").Value = -0.4306878306878307
$ws.Range("AQ4").Value = -32.56410256410256
